$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 22 through 29 results (rows 34-41), columns E/H/K/N/Q/T/W hold the
# raw scores for each of the 7 players; columns D/G/J/M/P/S/V already carry
# the VLOOKUP/RANK "placement points" formula and recalculate automatically.

$data = @{
    34 = @{ E=70;  H=40;  K=50;  N=80;  Q=100; T=0;   W=60 }
    35 = @{ E=80;  H=50;  K=60;  N=100; Q=70;  T=0;   W=40 }
    36 = @{ E=60;  H=50;  K=80;  N=70;  Q=100; T=40;  W=0  }
    37 = @{ E=80;  H=50;  K=40;  N=100; Q=0;   T=60;  W=70 }
    38 = @{ E=60;  H=80;  K=70;  N=40;  Q=100; T=40;  W=50 }
    39 = @{ E=70;  H=80;  K=40;  N=100; Q=0;   T=60;  W=50 }
    40 = @{ E=60;  H=100; K=70;  N=50;  Q=80;  T=40;  W=0  }
    41 = @{ E=0;   H=70;  K=80;  N=60;  Q=40;  T=100; W=50 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Row 38 has a tie between N38 and T38 (both 40), so RANK() returns the same
# rank for both and the VLOOKUP formula result is ambiguous; those two cells
# were overridden with the tie-break average value instead of the formula.
$ws.Range("M38").Value = -22.5
$ws.Range("S38").Value = -22.5

$excel.Calculate()
